$d = $word.ActiveDocument

# Word Find/Replace constants used below:
#   Wrap:      0 = wdFindStop, 1 = wdFindContinue
#   Replace:   1 = wdReplaceOne (replace only the first/next match), 2 = wdReplaceAll (document-wide!)
# NOTE: wdReplaceAll (2) replaces every match in the *entire* document regardless of
# which Range/Selection invoked Find, so whenever a single, specific occurrence must be
# targeted we use wdReplaceOne (1) on a Range anchored at/ before that occurrence.

# --------------------------------------------------------------------------
# Edit 1: Top summary box (Table 1, row 1, col 2) - TFS item title change.
#   "TFS 15833 - Warning Log Workflow Change"
#     -> "TFS 16370 - Director dashboard error due to follow-up code change"
#   Only the FIRST occurrence of "TFS 15833..." (the summary box) changes; the
#   identical text further down in the change-log table must stay untouched.
# --------------------------------------------------------------------------
$t1 = $word.ActiveDocument.Tables.Item(1)
$cell = $t1.Cell(1, 2)
$rng = $cell.Range
$rng.Find.Execute("5833", $false, $false, $false, $false, $false, $true, 0, $false, "6370", 1) | Out-Null

$t1 = $word.ActiveDocument.Tables.Item(1)
$cell = $t1.Cell(1, 2)
$rng = $cell.Range
$rng.Find.Execute("Warning Log Workflow Change", $false, $false, $false, $false, $false, $true, 0, $false, "Director dashboard error due to follow-up code change", 1) | Out-Null

# --------------------------------------------------------------------------
# Edit 2: Change-log table (Table 2) - append a new row documenting the change.
# --------------------------------------------------------------------------
$t2 = $word.ActiveDocument.Tables.Item(2)
$newRow = $t2.Rows.Add()

$t2 = $word.ActiveDocument.Tables.Item(2)
$newRow = $t2.Rows.Item($t2.Rows.Count)

$dateCell = $newRow.Cells.Item(1)
$dateCell.Range.Text = "01/08/2020"

$descCell = $newRow.Cells.Item(2)
$descCell.Range.Text = "TFS 16370 - Director dashboard error due to follow-up code change"

$authorCell = $newRow.Cells.Item(3)
$authorCell.Range.Text = "Lili Huang"

# --------------------------------------------------------------------------
# Edit 3: Purpose paragraph - collapse the TFS reference list down to the
# single new TFS number.
#   "... per TFS 14679, 15600, and 15653." -> "... per TFS 16370."
# --------------------------------------------------------------------------
$rng3 = $word.ActiveDocument.Content
$rng3.Find.Execute("14679, 15600, and 15653.", $false, $false, $false, $false, $false, $true, 0, $false, "16370.", 1) | Out-Null

# --------------------------------------------------------------------------
# Edit 4: Changeset number bump.
#   "Changeset 44175" -> "Changeset 44347"
# --------------------------------------------------------------------------
$rng4 = $word.ActiveDocument.Content
$rng4.Find.Execute("44175", $false, $false, $false, $false, $false, $true, 0, $false, "44347", 1) | Out-Null

# --------------------------------------------------------------------------
# Edit 5: Footer page-number cached field result on the 3rd section's footer
# (word/footer3.xml): "Page 3" -> "Page 2".
# --------------------------------------------------------------------------
$sec3 = $word.ActiveDocument.Sections.Item(3)
$ftr3 = $sec3.Footers.Item(1)
$ftrRng = $ftr3.Range
$ftrRng.Find.Execute("3", $false, $false, $false, $false, $false, $true, 0, $false, "2", 1) | Out-Null

Write-Host "Done."
